$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("A10").Value = 42612.883055555554
$ws.Range("B10").Value = 96
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = "Random"

# Row 11
$ws.Range("A11").Value = 42612.889594907407
$ws.Range("B11").Value = 7
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = "Random"
